# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) values were recalculated from the underlying
# box-score data and need to be rewritten with the new, correct strikeout
# counts for each of the 68 game rows (rows 2-69).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2 through 69, in order.
$kValues = @(
    1,0,2,1,1,2,2,2,1,2,
    0,1,2,2,1,2,1,1,1,0,
    0,1,2,3,1,1,1,0,2,0,
    3,1,1,1,1,1,1,1,3,0,
    1,2,2,1,1,1,0,2,2,0,
    0,2,2,1,0,3,0,2,1,2,
    3,1,3,2,1,2,0,1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}

Write-Host "Updated column G (K) for rows $startRow to $($startRow + $kValues.Length - 1)"
